$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- Helper: split a run at [splitStart, splitEnd) by toggling Bold on/off.
# This forces the underlying engine to materialize a standalone run for that
# sub-range without altering the run's effective formatting (rFonts survive
# intact, unlike Font.Name which drops w:cs).
function Split-Run($startOffset, $endOffset) {
    $rr = $d.Range($startOffset, $endOffset)
    $rr.Bold = 1
    $rr.Bold = 0
}

# =====================================================================
# 1) LinkedIn / tw / github line: merge the "https://...woods" run and the
#    "  * tw" run into a single run, keeping "103.github.io " a separate run.
# =====================================================================
$t = $d.Content.Text
$idxStart = $t.IndexOf("https://www.linkedin.com/in/tevariyae-woods")
$idxTwEnd = $t.IndexOf("103.github.io")
$r1 = $d.Range($idxStart, $idxTwEnd)

# Force a genuine text mutation (engine no-ops a Range.Text= whose value
# equals the current text), then restore the real target text.
$marker = [char]1 + "TEMP_MARKER" + [char]1
$r1.Text = $marker
$t2 = $d.Content.Text
$idxMarker = $t2.IndexOf($marker)
$r1b = $d.Range($idxMarker, $idxMarker + $marker.Length)
$r1b.Text = "https://www.linkedin.com/in/tevariyae-woods  " + $bullet + " tw"

# Re-establish "103.github.io " as its own run (the Text= above merges all
# runs in the paragraph into one).
$t3 = $d.Content.Text
$idxGh = $t3.IndexOf("103.github.io")
$idxGhEnd = $idxGh + "103.github.io ".Length
Split-Run $idxGh $idxGhEnd

# =====================================================================
# 2) Education bullet: "Bachelors of Science, Computer Science ..." ->
#    "Bachelor of Science" + ", Computer Science ..." (two runs).
# =====================================================================
$d.Content.Find.Execute("Bachelors of Science", $true, $false, $false, $false, $false, $true, 1, $false, "Bachelor of Science", 2) | Out-Null

$t4 = $d.Content.Text
$idxBach = $t4.IndexOf("Bachelor of Science")
$idxBachEnd = $idxBach + "Bachelor of Science".Length
Split-Run $idxBach $idxBachEnd

# =====================================================================
# 3) "This all time classicallowed a way ..." ->
#    "This " + "all-time" + " " + "classic allowed" + " a way ..."
#    (also clears the spellcheck proofErr markers around "classicallowed")
# =====================================================================
$t5 = $d.Content.Text
$phrase3 = $bullet + " This all time classicallowed"
$idxStart3 = $t5.IndexOf($phrase3)
# +1 consumes the trailing space / the proofErr "spellEnd" boundary so both
# proofErr elements fall inside the replaced range and get dropped.
$idxEnd3 = $idxStart3 + $phrase3.Length + 1
$r3 = $d.Range($idxStart3, $idxEnd3)
$r3.Text = $bullet + " This all-time classic allowed "

$b1s = $idxStart3 + 7
$b1e = $idxStart3 + 15
$b2e = $idxStart3 + 16
$b3e = $idxStart3 + 31
Split-Run $b1s $b1e
Split-Run $b1e $b2e
Split-Run $b2e $b3e

# =====================================================================
# 4) "... leading a group to achieving our goals; which was ..." ->
#    "... leading a group to " + "achieve" + " our " + "goals," + " which was ..."
# =====================================================================
$t6 = $d.Content.Text
$phrase4 = $bullet + " In this project I was one of the key coordinators leading a group to achieving our goals; which was to donate to the less fortunate."
$idxStart4 = $t6.IndexOf($phrase4)
$idxEnd4 = $idxStart4 + $phrase4.Length
$r4 = $d.Range($idxStart4, $idxEnd4)
$r4.Text = $bullet + " In this project I was one of the key coordinators leading a group to achieve our goals, which was to donate to the less fortunate."

$c1s = $idxStart4 + 71
$c1e = $idxStart4 + 78
$c2e = $idxStart4 + 83
$c3e = $idxStart4 + 89
Split-Run $c1s $c1e
Split-Run $c1e $c2e
Split-Run $c2e $c3e
